$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet from "Through 2022-02-11" to "Through 2022-02-12"
$ws.Name = "Through 2022-02-12"

# Update header label in I1 (shared string "2022 (through 02-11)" -> "2022 (through 02-12)")
$ws.Range("I1").Value = "2022 (through 02-12)"

# Update February total (I3): 52 -> 56
$ws.Range("I3").Value = 56

# Update overall Total row (I14): 213 -> 217
$ws.Range("I14").Value = 217
